{"js": "// Fix a spelling typo (\"palce\" -> \"place\") and flatten three paragraphs whose\n// runs were unnecessarily split (e.g. due to spell-check/grammar-check\n// proofing marks) back into single, unified runs with identical text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Map of a unique substring to find each target paragraph -> the full,\n// corrected text that paragraph's single run should contain afterwards.\nconst replacements = [\n  {\n    find: \"palce\",\n    text: \"-where does the level take place?\"\n  },\n  {\n    find: \"mains hobby\",\n    text: \"City hall is the mains hobby for getting missions, it is the first level, and the player should come back often. Also, the architecture is unique and grandiose.\"\n  },\n  {\n    find: \"The first level takes place\",\n    text: \"The first level takes place in the City Hall building. The player spawns in the lobby area, which features some armchairs and side tables made of wood and light fabric on both sides. There is also a reception in front of the player and an elevator behind the reception.\"\n  }\n];\n\nfor (const { find, text } of replacements) {\n  let target = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.indexOf(find) !== -1) {\n      target = paragraphs.items[i];\n      break;\n    }\n  }\n  if (target) {\n    // Replacing the whole paragraph range's text in one shot collapses all\n    // of its child runs (and any spell-check proofErr markers between them)\n    // into a single run carrying the corrected text.\n    target.getRange().insertText(text, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fix a spelling typo (\"palce\" -> \"place\") and flatten three paragraphs whose\n# runs were unnecessarily split (e.g. due to spell-check/grammar-check\n# proofing marks) back into single, unified runs with identical text.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ParagraphText($needle, $newText) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.Contains($needle)) {\n            # Building a fresh Range over the paragraph's [Start, End) span and\n            # assigning .Text replaces every run inside it (and drops any\n            # w:proofErr markers between them) with a single new run carrying\n            # the corrected text.\n            $r = $d.Range($p.Range.Start, $p.Range.End)\n            $r.Text = $newText\n            return\n        }\n    }\n}\n\nReplace-ParagraphText \"palce\" \"-where does the level take place?\"\nReplace-ParagraphText \"mains hobby\" \"City hall is the mains hobby for getting missions, it is the first level, and the player should come back often. Also, the architecture is unique and grandiose.\"\nReplace-ParagraphText \"The first level takes place\" \"The first level takes place in the City Hall building. The player spawns in the lobby area, which features some armchairs and side tables made of wood and light fabric on both sides. There is also a reception in front of the player and an elevator behind the reception.\"\n"}
